# Insert a new "Category" column before the existing "Spend" column (D),
# pushing Spend to column E, and populate it with category values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at D; existing column D (Spend) shifts to E.
$ws.Columns.Item(4).Insert()

# Header
$ws.Range("D1").Value = "Category"

# Category values for rows 2..21 (column D), keyed by row number.
$categories = @{
    2  = "D1"
    3  = "D2"
    4  = "D2"
    5  = "D3"
    6  = "D3"
    7  = "D1"
    8  = "D2"
    9  = "D2"
    10 = "D3"
    11 = "D3"
    12 = "D1"
    13 = "D1"
    14 = "D3"
    15 = "D5"
    16 = "D5"
    17 = "D1"
    18 = "D1"
    19 = "D2"
    20 = "D4"
    21 = "D4"
}

# Write the rows in the order that introduces each distinct category label
# for the first time as D1, D2, D3, D4, D5 (matching the original author's
# shared-string ordering), then the rest in natural order.
$writeOrder = @(2, 3, 5, 20, 15, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 21)

foreach ($row in $writeOrder) {
    $ws.Cells.Item($row, 4).Value = $categories[$row]
}

# Update the selection to match the saved view state.
[void]$ws.Range("D22").Select()
